$d = $word.ActiveDocument

# 1. "Pour les variables dependantes" paragraph: drop the trailing " 10" after
#    the final "<" so it reads "...dont le NOCom <" instead of "...dont le NOCom < 10"
$d.Content.Find.Execute(
    "la moyenne du DCP des classes dont le NOCom < 10",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "la moyenne du DCP des classes dont le NOCom <",
    2) | Out-Null

# 2. Append two new paragraphs at the end of the document (after the last
#    paragraph, before the section break), carrying on the same body text
#    formatting (rFonts cstheme=minorHAnsi) that the rest of the section uses.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs($lastIndex + 1)
$p1.Range.Text = "Le sondage a été fait à l’aide d’un programme appelé EtudeMetrique, qui a écrit les résultats dans un fichier appelé tp2#2.txt"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($lastIndex + 2)
$p2.Range.Text = "Selon tp2#2, la moyenne de la DCP des classes ayant un NOCom inférieure à 10 est 67.51805, alors que la moyenne de la DCP des classes ayant un NOCom supérieure à 10 est 48.18528, ce qui va directement à l’encontre de l’hypothèse posé. Il se pourrait que la plupart des commentaires est créé plus tôt durant le développement du code, au lieu d’être fait au fur et à mesure, ce qui donnerait donc des résultats tel que celui-ci."

Write-Output "done"
